# ---------------------------------------------------------------------------
# Edit summary (per the target diff):
#   1. The table on slide 6 (the "SOURCES OF FINANCE" table) switches its
#      table style from {5CAF36C9-2521-456C-A405-6568B6CB9D55} to
#      {1CF5A2FB-6CB3-4359-AFA6-6FC3AA923776}.
#   2. The presentation's theme (currently "Integral") is swapped for the
#      stock "Office Theme" palette - i.e. the design/theme that is applied
#      to the deck changes from the green "Integral" look to the default
#      blue/grey "Office Theme" look.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 --------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{1CF5A2FB-6CB3-4359-AFA6-6FC3AA923776}")

# --- 2. Swap the deck's theme colours from "Integral" to "Office Theme" ----
function HexToColorRef([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Standard "Office Theme" colour scheme (the 12 MsoThemeColorSchemeIndex
# slots, in order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToColorRef $officeThemeColors[$i - 1]
}
